$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 157 (shifts existing rows 157..243 down to 158..244)
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new data point
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44719
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100112032
$ws.Cells.Item(157, 7).Value = "Zapallo italiano"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 200
$ws.Cells.Item(157, 11).Value = 15000
$ws.Cells.Item(157, 12).Value = 15000
$ws.Cells.Item(157, 13).Value = 15000
$ws.Cells.Item(157, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(157, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(157, 16).Value = 300
$ws.Cells.Item(157, 17).Value = 50
$ws.Cells.Item(157, 18).Value = "Hortaliza"
